# Add a new "Behavioral traits" section (rows 29-35) to the explanatory
# variables table, mirroring the layout/style of the existing sections
# (e.g. "Life history traits" at A17, with its indented sub-items below).
#
# Values are written in the same order the source data was produced in
# (all of column A's new labels, then column B's new sources, then the
# new section header last) so the shared-string table fills in the same
# sequence as the target workbook. Formatting is copied (format-only
# paste) from existing cells that already carry the right look, which
# keeps the style table from growing with throwaway intermediate styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---- Column A labels (new rows 30-35), in first-use order -----------
$ws.Range("A31").Value = "Reproductive mode"
$ws.Range("A32").Value = "Reproductive guild 1"
$ws.Range("A33").Value = "Reproductive guild 2"
$ws.Range("A34").Value = "Spawning ground"
$ws.Range("A30").Value = "Migratory behavior"
$ws.Range("A35").Value = "Spawning frequency"

# ---- Column B sources (new rows 30-35), in first-use order ----------
$ws.Range("B31").Value = "Fishbase (finfish, 96%), SeaLifeBase (inverts, 95%)"
$ws.Range("B32").Value = "Fishbase (finfish, 91%), SeaLifeBase (inverts, 71%)"
$ws.Range("B33").Value = "Fishbase (finfish, 82%), SeaLifeBase (inverts, 67%)"
$ws.Range("B35").Value = "Fishbase (finfish, 55%), SeaLifeBase (inverts, 5%)"
$ws.Range("B34").Value = "Fishbase (finfish, 62%), SeaLifeBase (inverts, 0%)"
$ws.Range("B30").Value = "Fishbase (finfish, 69%), SeaLifeBase (inverts, 0%)"

# ---- New section header (row 29) -------------------------------------
$ws.Range("A29").Value = "Behavioral traits"

# --- Formatting ----------------------------------------------------
# Row 29: section header — same look as the other section headers
# (A3 "SST experience", A7 "Stock characteristics", A14 "Geography",
# A17 "Life history traits"): italic, left aligned, vertically centered.
# B29 already existed (empty) before this edit and is left untouched.
$ws.Range("A17").Copy() | Out-Null
$ws.Range("A29").PasteSpecial($xlPasteFormats) | Out-Null

# Row 30 "Migratory behavior" — vertically centered + indented sub-item
# (same look as the first sub-item under each section, e.g. A4).
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A30").PasteSpecial($xlPasteFormats) | Out-Null

# Rows 31-35: remaining sub-items, indented (same as A19-A28).
$ws.Range("A19").Copy() | Out-Null
"A31", "A32", "A33", "A34", "A35" | ForEach-Object {
    $ws.Range($_).PasteSpecial($xlPasteFormats) | Out-Null
}

# B30 & B31: plain left-aligned source cells (same as B19-B28).
$ws.Range("B19").Copy() | Out-Null
"B30", "B31" | ForEach-Object {
    $ws.Range($_).PasteSpecial($xlPasteFormats) | Out-Null
}

# B32-B35: left-aligned source cells with the darker/explicit-black font
# variant seen in the diff. Build it once on B32, then copy it onward.
$ws.Range("B32").Font.Color = 0
$ws.Range("B32").HorizontalAlignment = -4131
$ws.Range("B32").Copy() | Out-Null
"B33", "B34", "B35" | ForEach-Object {
    $ws.Range($_).PasteSpecial($xlPasteFormats) | Out-Null
}

$excel.CutCopyMode = $false

# Restore the selection to where the author left off.
$ws.Range("A19").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
